# Insert two new player rows (Jaime, Jorge Marcos) to keep the roster alphabetised,
# then refresh every statistic touched by the newly played game.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jugadores")

# Make room for the two new players (shifts rows 5.. and 7.. down by one each)
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(7).Insert()

# New player: Jaime
$ws.Range("A5").Value = "Jaime"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 5.28
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 2
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 33.3
$ws.Range("S5").Value = 50
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 50
$ws.Range("V5").Value = 35
$ws.Range("W5").Value = 1

# New player: Jorge Marcos
$ws.Range("A7").Value = "Jorge Marcos"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 4.97
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = 5
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 40
$ws.Range("S7").Value = 40
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 10
$ws.Range("W7").Value = 1

# Updated stats for existing players after the new game
# row 4
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 4
$ws.Range("G4").Value = 94.93000000000001
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 15
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 7
$ws.Range("N4").Value = 25
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 11
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 23.66
$ws.Range("S4").Value = 30
$ws.Range("T4").Value = 10
$ws.Range("U4").Value = 40
$ws.Range("W4").Value = 5

# row 8
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 43.32
$ws.Range("S8").Value = 43.32
$ws.Range("U8").Value = 20
$ws.Range("W8").Value = 5

# row 10
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 13
$ws.Range("G10").Value = 110.67
$ws.Range("H10").Value = 62
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 15
$ws.Range("L10").Value = 45
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 60
$ws.Range("O10").Value = 7
$ws.Range("P10").Value = 10
$ws.Range("Q10").Value = 1.4
$ws.Range("R10").Value = 35.38
$ws.Range("S10").Value = 29.66
$ws.Range("T10").Value = 35.12
$ws.Range("U10").Value = 60
$ws.Range("W10").Value = 5

# row 11
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 7
$ws.Range("E11").Value = 4
$ws.Range("G11").Value = 57.2
$ws.Range("H11").Value = 24
$ws.Range("J11").Value = 6
$ws.Range("K11").Value = 7
$ws.Range("L11").Value = 16
$ws.Range("M11").Value = 8
$ws.Range("N11").Value = 22
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 36
$ws.Range("S11").Value = 6.66
$ws.Range("T11").Value = 40
$ws.Range("U11").Value = 10
$ws.Range("W11").Value = 5

# row 12
$ws.Range("B12").Value = 15
$ws.Range("C12").Value = 13
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 7
$ws.Range("G12").Value = 81.5
$ws.Range("H12").Value = 43
$ws.Range("I12").Value = 17
$ws.Range("J12").Value = 24
$ws.Range("M12").Value = 17
$ws.Range("N12").Value = 27
$ws.Range("O12").Value = 9
$ws.Range("P12").Value = 12
$ws.Range("Q12").Value = 1.6
$ws.Range("R12").Value = 65.59999999999999
$ws.Range("S12").Value = 70.94
$ws.Range("U12").Value = 75
$ws.Range("W12").Value = 5

# row 13
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 13.05
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 3
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 3
$ws.Range("R13").Value = 25
$ws.Range("S13").Value = 25
$ws.Range("W13").Value = 2

# row 15
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 1
$ws.Range("G15").Value = 51.28
$ws.Range("H15").Value = 11
$ws.Range("J15").Value = 8
$ws.Range("N15").Value = 9
$ws.Range("O15").Value = 3
$ws.Range("P15").Value = 7
$ws.Range("R15").Value = 20
$ws.Range("S15").Value = 20
$ws.Range("U15").Value = 12.5
$ws.Range("W15").Value = 4

# row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 4
$ws.Range("G16").Value = 59.05
$ws.Range("H16").Value = 52
$ws.Range("J16").Value = 6
$ws.Range("K16").Value = 13
$ws.Range("L16").Value = 36
$ws.Range("M16").Value = 16
$ws.Range("N16").Value = 42
$ws.Range("Q16").Value = 2.5
$ws.Range("R16").Value = 39.35
$ws.Range("S16").Value = 58.32
$ws.Range("T16").Value = 38.12
$ws.Range("U16").Value = 35.42
$ws.Range("W16").Value = 4

# row 17
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 8
$ws.Range("G17").Value = 56.62
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 7
$ws.Range("L17").Value = 4
$ws.Range("M17").Value = 4
$ws.Range("N17").Value = 11
$ws.Range("Q17").Value = 2.33
$ws.Range("R17").Value = 30
$ws.Range("S17").Value = 55.57
$ws.Range("U17").Value = 25
$ws.Range("W17").Value = 3

# row 18
$ws.Range("B18").Value = 27
$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 109.95
$ws.Range("H18").Value = 51
$ws.Range("I18").Value = 19
$ws.Range("J18").Value = 36
$ws.Range("M18").Value = 19
$ws.Range("N18").Value = 37
$ws.Range("O18").Value = 13
$ws.Range("P18").Value = 24
$ws.Range("Q18").Value = 1.4
$ws.Range("R18").Value = 53
$ws.Range("S18").Value = 54.34
$ws.Range("U18").Value = 60.6
$ws.Range("W18").Value = 5

# Team averages (Stats_Rival) recomputed for the new game
$ws2 = $wb.Worksheets.Item("Stats_Rival")
$ws2.Range("A2").Value = 74.8
$ws2.Range("B2").Value = 53.4
$ws2.Range("C2").Value = 21.8
$ws2.Range("D2").Value = 29.4
$ws2.Range("E2").Value = 51.2
$ws2.Range("F2").Value = 7.8
$ws2.Range("G2").Value = 22.6
$ws2.Range("H2").Value = 30.4
